$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column values that look numeric (e.g. "165.26") are written as literal
# text, matching the existing inline-string/text storage used throughout column D,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.145.65"
$ws.Range("E2").Value = "  -4.31%  "

$ws.Range("D3").Value = "2.490.08"
$ws.Range("E3").Value = "  -5.57%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "572.77"
$ws.Range("E5").Value = "  -2.78%  "

$ws.Range("D6").Value = "165.26"
$ws.Range("E6").Value = "  -4.95%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "0.512"
$ws.Range("E8").Value = "  -1.23%  "

$ws.Range("D9").Value = "2.487.35"
$ws.Range("E9").Value = "  -5.63%  "

$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  -10.18%  "

$ws.Range("D11").Value = "0.167"
$ws.Range("E11").Value = "  -1.26%  "

$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  -4.42%  "

$ws.Range("D13").Value = "4.81"
$ws.Range("E13").Value = "  -2.56%  "

$ws.Range("D14").Value = "2.950.60"
$ws.Range("E14").Value = "  -5.79%  "

$ws.Range("D15").Value = "69.088.06"
$ws.Range("E15").Value = "  -4.21%  "

$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  -6.71%  "

$ws.Range("D17").Value = "24.46"
$ws.Range("E17").Value = "  -5.33%  "

$ws.Range("D18").Value = "2.494.84"
$ws.Range("E18").Value = "  -5.62%  "

$ws.Range("D19").Value = "11.29"
$ws.Range("E19").Value = "  -6.85%  "

$ws.Range("D20").Value = "7.70"
$ws.Range("E20").Value = "  -3.12%  "

$ws.Range("D21").Value = "344.78"
$ws.Range("E21").Value = "  -7.26%  "

$ws.Range("D22").Value = "3.90"
$ws.Range("E22").Value = "  -5.74%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").Value = "1.93"
$ws.Range("E24").Value = "  -5.39%  "

$ws.Range("D25").Value = "67.98"
$ws.Range("E25").Value = "  -4.18%  "

$ws.Range("D26").Value = "3.93"
$ws.Range("E26").Value = "  -7.27%  "

$ws.Range("D27").Value = "8.83"
$ws.Range("E27").Value = "  -8.16%  "

$ws.Range("D28").Value = "2.637.85"
$ws.Range("E28").Value = "  -4.90%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("D30").Value = "0.0₃0887"
$ws.Range("E30").Value = "  -6.71%  "

$ws.Range("D31").Value = "7.77"
$ws.Range("E31").Value = "  -2.52%  "

$ws.Range("D32").Value = "464.43"
$ws.Range("E32").Value = "  -6.28%  "

$ws.Range("D33").Value = "1.25"
$ws.Range("E33").Value = "  -2.16%  "

$ws.Range("D34").Value = "1.73"
$ws.Range("E34").Value = "  -3.83%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("E36").Value = "  +1.66%  "

$ws.Range("D37").Value = "152.46"
$ws.Range("E37").Value = "  -5.52%  "

$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("D39").Value = "18.28"
$ws.Range("E39").Value = "  -4.95%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").Value = "4.69"
$ws.Range("E41").Value = "  -4.24%  "

$ws.Range("D42").Value = "0.313"
$ws.Range("E42").Value = "  -3.96%  "

$ws.Range("D43").Value = "1.57"
$ws.Range("E43").Value = "  -9.40%  "

$ws.Range("D44").Value = "1.15"
$ws.Range("E44").Value = "  -14.72%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "38.07"
$ws.Range("E45").Value = "  -2.55%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.28"
$ws.Range("E46").Value = "  -10.76%  "

$ws.Range("D47").Value = "141.73"
$ws.Range("E47").Value = "  -7.10%  "

$ws.Range("D48").Value = "0.521"
$ws.Range("E48").Value = "  -4.54%  "

$ws.Range("D49").Value = "3.47"
$ws.Range("E49").Value = "  -4.77%  "

$ws.Range("E50").Value = "  -5.60%  "

$ws.Range("D51").Value = "0.0730"
$ws.Range("E51").Value = "  -2.25%  "

# Restore default (unstyled) formatting for column D so only cell contents change.
$ws.Range("D2:D51").Style = "Normal"
